$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting (bold, bordered, centered) used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill season record columns for every data row (2-46)
$ws.Range("AD2:AD46").Value = 83
$ws.Range("AE2:AE46").Value = 79
$ws.Range("AF2:AF46").Value = 0
